$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are treated as text so numeric-looking strings
# (e.g. "657.65", "1.00") are preserved exactly, matching the source data.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '96.969.50'
$ws.Range('D3').Value = '3.344.15'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '250.75'
$ws.Range('E5').Value = '  -1.84%  '
$ws.Range('D6').Value = '657.65'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').Value = '1.43'
$ws.Range('E7').Value = '  -4.48%  '
$ws.Range('D8').Value = '0.425'
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '1.02'
$ws.Range('E10').Value = '  -5.51%  '
$ws.Range('D11').Value = '3.341.55'
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('D12').Value = '0.207'
$ws.Range('E12').Value = '  -2.73%  '
$ws.Range('D13').Value = '41.04'
$ws.Range('E13').Value = '  -2.66%  '
$ws.Range('D14').Value = '96.722.05'
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').Value = '6.11'
$ws.Range('E15').Value = '  -3.76%  '
$ws.Range('E16').Value = '  -2.59%  '
$ws.Range('D17').Value = '3.967.53'
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D18').Value = '8.83'
$ws.Range('E18').Value = '  +2.26%  '
$ws.Range('D19').Value = '3.320.54'
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').Value = '0.559'
$ws.Range('E20').Value = '  +11.09%  '
$ws.Range('D21').Value = '17.65'
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').Value = '10.73'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '510.19'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('E24').Value = '  -3.80%  '
$ws.Range('E25').Value = '  -3.72%  '
$ws.Range('D26').Value = '6.61'
$ws.Range('E26').Value = '  +6.05%  '
$ws.Range('D27').Value = '96.86'
$ws.Range('E27').Value = '  -2.57%  '
$ws.Range('D28').Value = '12.17'
$ws.Range('E28').Value = '  -5.42%  '
$ws.Range('E29').Value = '  -4.62%  '
$ws.Range('D30').Value = '11.49'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('E32').Value = '  -6.85%  '
$ws.Range('E33').Value = '  +10.81%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '0.558'
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('D36').Value = '28.47'
$ws.Range('E36').Value = '  -4.66%  '
$ws.Range('E37').Value = '  +4.12%  '
$ws.Range('D38').Value = '7.84'
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.152'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').Value = '508.17'
$ws.Range('E41').Value = '  -3.35%  '
$ws.Range('D42').Value = '0.0438'
$ws.Range('E42').Value = '  +4.04%  '
$ws.Range('D43').Value = '24.38'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('D44').Value = '0.840'
$ws.Range('E44').Value = '  -3.56%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value = '5.67'
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D46').Value = '3.65'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('D47').Value = '1.68'
$ws.Range('E47').Value = '  +5.07%  '
$ws.Range('D48').Value = '8.59'
$ws.Range('E48').Value = '  +3.16%  '
$ws.Range('D49').Value = '54.71'
$ws.Range('E49').Value = '  +6.42%  '
$ws.Range('E50').Value = '  -5.90%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '2.00'
$ws.Range('E51').Value = '  -3.84%  '

# Restore default (Normal) style so no stray number-format styling remains
# on the edited cells, matching the original workbook formatting.
$ws.Range("B2:E51").Style = "Normal"
